# Auto-generated Excel COM-interop script updating crypto price/volume cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.603.07"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "1.698.78"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.41"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3728"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.63"
$ws.Range("E8").Value = "  +2.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3421"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.175"
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.74"
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.884"
$ws.Range("D16").Value = "1.702.04"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001116"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06658"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "83.00"
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.00"
$ws.Range("E21").Value = "  +2.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.312"
$ws.Range("E22").Value = "  +2.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.12"
$ws.Range("E23").Value = "  +8.91%  "
$ws.Range("D24").Value = "24.563.42"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.430"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("E27").Value = "  +2.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.06"
$ws.Range("E28").Value = "  -1.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "130.52"
$ws.Range("E29").Value = "  +2.81%  "
$ws.Range("D30").Value = "1.889.14"
$ws.Range("E30").Value = "  +1.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.167"
$ws.Range("E31").Value = "  +17.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.667"
$ws.Range("E32").Value = "  +3.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.211"
$ws.Range("E33").Value = "  +3.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08744"
$ws.Range("E35").Value = "  -1.10%  "
$ws.Range("E36").Value = "  +7.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.470"
$ws.Range("E37").Value = "  +1.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06484"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.875"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02351"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2173"
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.270"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6360"
$ws.Range("E43").Value = "  +2.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9996"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.82"
$ws.Range("E45").Value = "  +4.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6037"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("E48").Value = "  +2.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "128.35"
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07226"
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.50"
$ws.Range("E51").Value = "  +2.50%  "
